$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings that must stay as TEXT
# (preserve exact formatting, e.g. trailing zeros / multi-dot thousands separators).
# Pre-format as Text so Excel does not silently convert them to numbers.
$textCells = @("D5", "D6", "D8", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D28", "D30", "D31", "D33", "D36", "D39", "D40", "D42", "D44", "D46", "D47", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values, in sheet order.
$ws.Range("D2").Value = '67.940.90'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '3.270.09'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '186.18'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").Value = '581.33'
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("D9").Value = '3.265.85'
$ws.Range("E9").Value = '  -1.19%  '
$ws.Range("E10").Value = '  -3.89%  '
$ws.Range("E11").Value = '  -2.20%  '
$ws.Range("E12").Value = '  -2.84%  '
$ws.Range("D13").Value = '3.836.69'
$ws.Range("E13").Value = '  -1.13%  '
$ws.Range("D14").Value = '0.137'
$ws.Range("D15").Value = '27.57'
$ws.Range("E15").Value = '  -5.16%  '
$ws.Range("D16").Value = '67.949.35'
$ws.Range("E16").Value = '  -1.31%  '
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("D18").Value = '3.263.39'
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("D19").Value = '5.76'
$ws.Range("E19").Value = '  -2.70%  '
$ws.Range("D20").Value = '13.56'
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").Value = '399.19'
$ws.Range("E21").Value = '  +3.36%  '
$ws.Range("D22").Value = '7.64'
$ws.Range("E22").Value = '  -2.44%  '
$ws.Range("D23").Value = '71.55'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("E26").Value = '  -4.23%  '
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("D28").Value = '9.52'
$ws.Range("E28").Value = '  -3.82%  '
$ws.Range("E29").Value = '  +0.47%  '
$ws.Range("D30").Value = '1.96'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("D31").Value = '22.71'
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("E32").Value = '  -6.77%  '
$ws.Range("D33").Value = '6.96'
$ws.Range("E33").Value = '  -4.17%  '
$ws.Range("E34").Value = '  -6.06%  '
$ws.Range("D36").Value = '163.51'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  -5.40%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").Value = '26.76'
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("D40").Value = '0.810'
$ws.Range("E40").Value = '  -3.76%  '
$ws.Range("E41").Value = '  -2.18%  '
$ws.Range("D42").Value = '6.42'
$ws.Range("E42").Value = '  -5.52%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.673.10'
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").Value = '0.0686'
$ws.Range("E44").Value = '  -1.79%  '
$ws.Range("E45").Value = '  -1.70%  '
$ws.Range("D46").Value = '2.44'
$ws.Range("E46").Value = '  -8.51%  '
$ws.Range("D47").Value = '24.83'
$ws.Range("E47").Value = '  -4.14%  '
$ws.Range("D48").Value = '333.52'
$ws.Range("E48").Value = '  -2.73%  '
$ws.Range("E49").Value = '  -3.41%  '
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("E51").Value = '  -1.66%  '
